{"js": "// Fix the title of the data analysis report:\n//   \"Data Analysis \u2013 Matplotlib Challenge\"  ->  \"Data Analysis \u2013 API Challenge\"\n// and move the \"_GoBack\" bookmark from the end of the document (after the\n// third bullet) to the title line, right after the new \"API\" text.\n\n// 1) Remove the stray \"_GoBack\" bookmark currently sitting at the end of the\n//    third paragraph (it will be re-inserted at the title).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Replace \"Matplotlib\" with \"API\" in the title paragraph.\nconst titleMatches = context.document.body.search(\"Matplotlib\", { matchCase: true });\ntitleMatches.load(\"text\");\nawait context.sync();\n\nif (titleMatches.items.length > 0) {\n  titleMatches.items[0].insertText(\"API\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Re-locate the word \"API\" we just inserted and drop the \"_GoBack\"\n//    bookmark immediately after it (i.e. between \"API\" and \" Challenge\").\nconst apiMatches = context.document.body.search(\"API\", { matchCase: true });\napiMatches.load(\"text\");\nawait context.sync();\n\nif (apiMatches.items.length > 0) {\n  const afterApi = apiMatches.items[0].getRange(\"After\");\n  afterApi.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Fix the title of the data analysis report:\n#   \"Data Analysis - Matplotlib Challenge\"  ->  \"Data Analysis - API Challenge\"\n# and move the \"_GoBack\" bookmark from the end of the document (after the\n# third bullet) to the title line, right after the new \"API\" text.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the stray \"_GoBack\" bookmark currently sitting at the end of the\n#    third paragraph (it will be re-inserted at the title below).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Replace \"Matplotlib\" with \"API\" in the title paragraph.\n$findRange = $d.Content\n$findRange.Find.Execute(\"Matplotlib\", $false, $false, $false, $false, $false, $true, 1, $false, \"API\", 2)\n\n# 3) Re-locate the word \"API\" we just inserted, collapse the range to its\n#    end, and drop the \"_GoBack\" bookmark there (i.e. between \"API\" and\n#    \" Challenge\").\n$apiRange = $d.Content\n$found = $apiRange.Find.Execute(\"API\")\nif ($found) {\n    $apiRange.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $apiRange)\n}\n"}
